# Apply the "new changes with added test cases" commit:
#  - D12 (RESPONSE NAME = utter_latest_publications) TEXT value gets
#    target="_blank" added to each <a> tag.
#  - Row 12 grows from 90pt to 135pt to fit the now-longer wrapped text.
#  - Selection moves from D4 to D12 (the cell that was just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "Read more on the latest in the world of business and technology: \n`n<a href=""https://www.marlabs.com/digital-victories/"" target=""_blank"">Case studies</a> | <a href=""https://www.marlabs.com/now-and-next/"" target=""_blank"">Whitepapers</a> | <a href=""https://www.marlabs.com/now-and-next/"" target=""_blank"">Podcasts</a> | <a href=""https://www.marlabs.com/now-and-next/"" target=""_blank"">Videos</a>"

$ws.Range("D12").Value = $newText

$ws.Rows.Item(12).RowHeight = 135

$ws.Range("D12").Select()
